# Build 70 intermediate push: add USR_ALERT_1/2/3 and THOOK columns to the
# Tabelle2 (sheet2) merge/lookup table, inserted right before the
# END_OF_COL / Title columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Insert 4 new blank columns at FC:FF - this pushes the existing
# END_OF_COL (old FC) and Title (old FD) columns to FG and FH,
# carrying their formulas/styles/widths along automatically.
$ws.Range("FC:FF").Insert()

# ---- Header row (row 1) ----
$ws.Range("FC1").Value = "USR_ALERT_1"
$ws.Range("FD1").Value = "USR_ALERT_2"
$ws.Range("FE1").Value = "USR_ALERT_3"
$ws.Range("FF1").Value = "THOOK"

# ---- Data rows (rows 2-41) ----
$ws.Range("FC2:FF41").Value = "|"

# Row 5 (Boeing F/A 18E Super Hornet) has a numeric marker in THOOK (FF)
$ws.Range("FF5").Value = 1

# Match the author's last-used selection from the commit
$ws.Range("FH9").Select()
